# Updates the LR-pair table (rows 2-11) with newly computed TPM-based values.
# The sending/target cluster pairing now also includes "MuSCs" as a target
# cluster (previously only "FAPs"), doubling the number of data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds the full row of values for columns A..T, in order:
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity (avg), Ligand derived specificity (total),
# Receptor-expressing cells, Receptor detection rate, Receptor average expression value,
# Receptor total expression value, Receptor derived specificity (avg), Receptor derived specificity (total),
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity
$rows = @(
    @("ECs", "Tnc", "Cntn1", "FAPs", 3, 1, 7.759382333333334, 23.278147, 0.03884312320086166, 0.04014625174390325, 2, 0.6666666666666666, 0.07031566666666667, 0.210947, 0.3500709860127268, 0.446885632088942, 0.5456061416898889, 4.910455275209, 0.01359785043873947, 0.01794078308657599),
    @("ECs", "Tnc", "Cntn1", "MuSCs", 3, 1, 7.759382333333334, 23.278147, 0.03884312320086166, 0.04014625174390325, 2, 1, 0.1305455, 0.261091, 0.6499290139872732, 0.5531143679110581, 1.012952446396167, 6.077714678377, 0.02524527276212219, 0.02220546865732726),
    @("FAPs", "Tnc", "Cntn1", "FAPs", 3, 1, 172.4328083333333, 517.298425, 0.8631909770948131, 0.8921497401307179, 2, 0.6666666666666666, 0.07031566666666667, 0.210947, 0.3500709860127268, 0.446885632088942, 12.12472787316389, 109.122550858475, 0.3021781164688703, 0.3986889005363012),
    @("FAPs", "Tnc", "Cntn1", "MuSCs", 3, 1, 172.4328083333333, 517.298425, 0.8631909770948131, 0.8921497401307179, 2, 1, 0.1305455, 0.261091, 0.6499290139872732, 0.5531143679110581, 22.51032718027917, 135.061963081675, 0.5610128606259428, 0.4934608395944167),
    @("Inflammatory-Mac", "Tnc", "Cntn1", "FAPs", 2, 0.6666666666666666, 0.05240566666666666, 0.157217, 0.0002623404388789996, 0.0002711415672571033, 2, 0.6666666666666666, 0.07031566666666667, 0.210947, 0.3500709860127268, 0.446885632088942, 0.003684939388777777, 0.033164454499, [double]"9.183777610938288e-05", 0.000121169270669277),
    @("Inflammatory-Mac", "Tnc", "Cntn1", "MuSCs", 2, 0.6666666666666666, 0.05240566666666666, 0.157217, 0.0002623404388789996, 0.0002711415672571033, 2, 1, 0.1305455, 0.261091, 0.6499290139872732, 0.5531143679110581, 0.006841323957833333, 0.041047943747, 0.0001705026627696167, 0.0001499722965878263),
    @("MuSCs", "Tnc", "Cntn1", "FAPs", 2, 1, 19.452549, 38.905098, 0.09737859599105524, 0.06709700125311635, 2, 0.6666666666666666, 0.07031566666666667, 0.210947, 0.3500709860127268, 0.446885632088942, 1.367818951301, 8.206913707805999, 0.03408942111512368, 0.02998468581627143),
    @("MuSCs", "Tnc", "Cntn1", "MuSCs", 2, 1, 19.452549, 38.905098, 0.09737859599105524, 0.06709700125311635, 2, 1, 0.1305455, 0.261091, 0.6499290139872732, 0.5531143679110581, 2.5394427354795, 10.157770941918, 0.06328917487593157, 0.03711231543684492),
    @("Resolving-Mac", "Tnc", "Cntn1", "FAPs", 1, 0.3333333333333333, 0.06491533333333334, 0.194746, 0.0003249632743909987, 0.0003358653050055137, 2, 0.6666666666666666, 0.07031566666666667, 0.210947, 0.3500709860127268, 0.446885632088942, 0.004564564940222222, 0.041081084462, 0.0001137602138839812, 0.0001500933791241343),
    @("Resolving-Mac", "Tnc", "Cntn1", "MuSCs", 1, 0.3333333333333333, 0.06491533333333334, 0.194746, 0.0003249632743909987, 0.0003358653050055137, 2, 1, 0.1305455, 0.261091, 0.6499290139872732, 0.5531143679110581, 0.008474404647666668, 0.050846427886, 0.0002112030605070175, 0.0001857719258813794)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

Write-Host "Updated rows 2-11 with recomputed TPM-based LR-pair values"
